# totul functioneaza. trebuie sa aranjez codul si sa schimb culorile
#
# - recompute the "Realized"/"Progress"/"Label" roll-up columns (E/F/G) for the
#   existing goals/tasks after task list changes
# - delete the two tasks that are gone ("aba" and "FINISHHHHH")
# - append three brand-new tasks ("aNewTaskForANewGoal", "anotherTask", "nTask")
# - move the selection to E4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- update the roll-up columns on the rows that survive untouched ---------
# Goal2
$ws.Range("E3").Value = "1|1"
$ws.Range("G3").Value = "Just Started"
# Task2.2
$ws.Range("F4").Value = "2%"
$ws.Range("G4").Value = "Just Started"
# Task2.3
$ws.Range("E5").Value = "1|1"
$ws.Range("G5").Value = "Just Started"

# --- remove the two finished/deleted tasks ("aba" row 6, "FINISHHHHH" row 7) ---
$ws.Rows("6:7").Delete()

# after the delete, the old rows 8-12 are now rows 6-10
# newGoal
$ws.Range("E6").Value = "1|1"
# aTaskBecauseTheOtherWasDeleted
$ws.Range("E7").Value = "1|3"
# uguigiuhiuh
$ws.Range("E8").Value = "1|1"
# guyguih9o
$ws.Range("E9").Value = "1|1"
$ws.Range("F9").Value = "15%"
# brandNewGoal
$ws.Range("E10").Value = "1|1"

# --- append the new tasks --------------------------------------------------
$ws.Range("A11").Value = "Task"
$ws.Range("B11").Value = "aNewTaskForANewGoal"
$ws.Range("C11").Value = 43874.49643962963
$ws.Range("C11").NumberFormat = "dd/MM/yyyy"
$ws.Range("D11").Value = "2|3"
$ws.Range("E11").Value = "0|0"
$ws.Range("F11").Value = "25%"
$ws.Range("G11").Value = "Something Done"

$ws.Range("A12").Value = "Task"
$ws.Range("B12").Value = "anotherTask"
$ws.Range("C12").Value = 44253.51310829861
$ws.Range("C12").NumberFormat = "dd/MM/yyyy"
$ws.Range("D12").Value = "4|6"
$ws.Range("E12").Value = "0|0"
$ws.Range("F12").Value = "0%"
$ws.Range("G12").Value = "Just Started"

$ws.Range("A13").Value = "Task"
$ws.Range("B13").Value = "nTask"
$ws.Range("C13").Value = 44233.51383005787
$ws.Range("C13").NumberFormat = "dd/MM/yyyy"
$ws.Range("D13").Value = "5|6"
$ws.Range("E13").Value = "0|0"
$ws.Range("F13").Value = "0%"
$ws.Range("G13").Value = "Just Started"

# --- move the selection, like the author did while poking at the new rows --
$ws.Range("E4").Select()
